# Insert a new date column ("09-dec") into the "Prix Spot" sheet just
# before the existing "01-oct." column (column EF), shifting every
# column from EF onward one place to the right (EF -> EG, ..., FJ -> FK).
# The new column gets the header "09-dec" in row 1 and "-" placeholders
# for the data rows (2-25), matching the existing "no data" convention
# used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Shift columns EF:FJ -> EG:FK, leaving a blank column at EF.
$ws.Columns("EF:EF").Insert()

# Populate the newly inserted column.
$ws.Range("EF1").Value = "09-dec"
$ws.Range("EF2:EF25").Value = "-"
